$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Value = "277.91" }
    @{ Addr = "E2"; Value = "1.57%" }
    @{ Addr = "D3"; Value = "27.35" }
    @{ Addr = "E3"; Value = "1.96%" }
    @{ Addr = "D4"; Value = "4.900" }
    @{ Addr = "E4"; Value = "-0.43%" }
    @{ Addr = "D5"; Value = "0.06370" }
    @{ Addr = "E5"; Value = "0.64%" }
    @{ Addr = "D6"; Value = "6.976" }
    @{ Addr = "E6"; Value = "0.51%" }
    @{ Addr = "D7"; Value = "1.266" }
    @{ Addr = "E7"; Value = "-7.30%" }
    @{ Addr = "D8"; Value = "0.8814" }
    @{ Addr = "E8"; Value = "-0.75%" }
    @{ Addr = "D9"; Value = "0.1526" }
    @{ Addr = "E9"; Value = "3.65%" }
    @{ Addr = "D10"; Value = "0.05108" }
    @{ Addr = "E10"; Value = "-0.49%" }
    @{ Addr = "D11"; Value = "0.07595" }
    @{ Addr = "E11"; Value = "3.37%" }
    @{ Addr = "D12"; Value = "0.02966" }
    @{ Addr = "E12"; Value = "-6.64%" }
    @{ Addr = "D13"; Value = "0.09007" }
    @{ Addr = "E13"; Value = "-0.68%" }
    @{ Addr = "D14"; Value = "0.001567" }
    @{ Addr = "E14"; Value = "0.09%" }
    @{ Addr = "D15"; Value = "0.0006406" }
    @{ Addr = "E15"; Value = "1.46%" }
    @{ Addr = "D16"; Value = "0.005966" }
    @{ Addr = "E16"; Value = "-1.22%" }
    @{ Addr = "D17"; Value = "3.462" }
    @{ Addr = "E17"; Value = "-0.43%" }
    @{ Addr = "D18"; Value = "3.310" }
    @{ Addr = "E18"; Value = "-1.33%" }
    @{ Addr = "D19"; Value = "2.272" }
    @{ Addr = "E19"; Value = "-0.40%" }
    @{ Addr = "E20"; Value = "0.53%" }
    @{ Addr = "E21"; Value = "0.32%" }
    @{ Addr = "E22"; Value = "-0.16%" }
    @{ Addr = "D23"; Value = "0.04420" }
    @{ Addr = "E23"; Value = "1.56%" }
    @{ Addr = "D24"; Value = "0.001171" }
    @{ Addr = "E24"; Value = "-0.87%" }
    @{ Addr = "D25"; Value = "0.003872" }
    @{ Addr = "E25"; Value = "6.26%" }
    @{ Addr = "D26"; Value = "0.0001199" }
    @{ Addr = "E26"; Value = "-0.49%" }
    @{ Addr = "E27"; Value = "-0.48%" }
    @{ Addr = "D40"; Value = "0.04147" }
    @{ Addr = "E40"; Value = "2.89%" }
    @{ Addr = "D41"; Value = "0.006812" }
    @{ Addr = "E41"; Value = "2.90%" }
    @{ Addr = "D42"; Value = "0.1178" }
    @{ Addr = "E42"; Value = "1.05%" }
    @{ Addr = "D43"; Value = "0.002069" }
    @{ Addr = "E43"; Value = "-12.73%" }
    @{ Addr = "D44"; Value = "0.01123" }
    @{ Addr = "E44"; Value = "-11.01%" }
    @{ Addr = "D45"; Value = "0.00005165" }
    @{ Addr = "E45"; Value = "-1.84%" }
    @{ Addr = "D47"; Value = "0.02023" }
    @{ Addr = "E47"; Value = "-5.00%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
